$d = $word.ActiveDocument

$d.Content.Find.Execute("78×47=", $true, $false, $false, $false, $false, $true, 1, $false, "97×43=", 2) | Out-Null
$d.Content.Find.Execute("44×29=", $true, $false, $false, $false, $false, $true, 1, $false, "66×48=", 2) | Out-Null
$d.Content.Find.Execute("77×67=", $true, $false, $false, $false, $false, $true, 1, $false, "69×96=", 2) | Out-Null
$d.Content.Find.Execute("17×11=", $true, $false, $false, $false, $false, $true, 1, $false, "54×33=", 2) | Out-Null
$d.Content.Find.Execute("14×90=", $true, $false, $false, $false, $false, $true, 1, $false, "86×58=", 2) | Out-Null
$d.Content.Find.Execute("37×27=", $true, $false, $false, $false, $false, $true, 1, $false, "68×89=", 2) | Out-Null
$d.Content.Find.Execute("97×42=", $true, $false, $false, $false, $false, $true, 1, $false, "33×72=", 2) | Out-Null
$d.Content.Find.Execute("98×69=", $true, $false, $false, $false, $false, $true, 1, $false, "67×38=", 2) | Out-Null
$d.Content.Find.Execute("41×86=", $true, $false, $false, $false, $false, $true, 1, $false, "78×45=", 2) | Out-Null
$d.Content.Find.Execute("58×13=", $true, $false, $false, $false, $false, $true, 1, $false, "69×17=", 2) | Out-Null
$d.Content.Find.Execute("28×85=", $true, $false, $false, $false, $false, $true, 1, $false, "72×82=", 2) | Out-Null
$d.Content.Find.Execute("66×92=", $true, $false, $false, $false, $false, $true, 1, $false, "17×87=", 2) | Out-Null
$d.Content.Find.Execute("32×84=", $true, $false, $false, $false, $false, $true, 1, $false, "94×68=", 2) | Out-Null
$d.Content.Find.Execute("47×95=", $true, $false, $false, $false, $false, $true, 1, $false, "83×44=", 2) | Out-Null
$d.Content.Find.Execute("77×96=", $true, $false, $false, $false, $false, $true, 1, $false, "57×16=", 2) | Out-Null
$d.Content.Find.Execute("86×67=", $true, $false, $false, $false, $false, $true, 1, $false, "94×43=", 2) | Out-Null
$d.Content.Find.Execute("94×24=", $true, $false, $false, $false, $false, $true, 1, $false, "59×39=", 2) | Out-Null
$d.Content.Find.Execute("47×25=", $true, $false, $false, $false, $false, $true, 1, $false, "84×61=", 2) | Out-Null
$d.Content.Find.Execute("38×81=", $true, $false, $false, $false, $false, $true, 1, $false, "83×83=", 2) | Out-Null
$d.Content.Find.Execute("32×59=", $true, $false, $false, $false, $false, $true, 1, $false, "69×63=", 2) | Out-Null
$d.Content.Find.Execute("75×47=", $true, $false, $false, $false, $false, $true, 1, $false, "14×15=", 2) | Out-Null
$d.Content.Find.Execute("80×58=", $true, $false, $false, $false, $false, $true, 1, $false, "75×44=", 2) | Out-Null
$d.Content.Find.Execute("61×88=", $true, $false, $false, $false, $false, $true, 1, $false, "62×83=", 2) | Out-Null
$d.Content.Find.Execute("39×78=", $true, $false, $false, $false, $false, $true, 1, $false, "34×79=", 2) | Out-Null
$d.Content.Find.Execute("56×68=", $true, $false, $false, $false, $false, $true, 1, $false, "39×83=", 2) | Out-Null
